# Apply the "added some more codes : Manish" edit to the workbook.
#
# Summary of changes:
#  - Sheet2: C2 "Motorized Card Reader" -> "manish", C3 "Motorized Card Reader" -> "ashish"
#  - A new Sheet3 is added after Sheet2 containing a small device-status table
#    for 4 recycler devices (IDC01S1..IDC01S4)
#  - Selection/active-sheet bookkeeping is updated so Sheet3 ends up the
#    active (tab-selected) sheet, Sheet1 ends up with the whole sheet
#    selected and Sheet2 ends up with C1 selected.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# 1. Sheet2 data tweak: rename the two "Motorized Card Reader" device
#    names to the two engineers who contributed the extra codes.
# ---------------------------------------------------------------------
$ws2.Range("C2").Value = "manish"
$ws2.Range("C3").Value = "ashish"

# ---------------------------------------------------------------------
# 2. Insert the new Sheet3 right after Sheet2.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Sheet3"

# Header row: statusCode / IDC01S1 / IDC01S2 / IDC01S3 / IDC01S4
$ws3.Range("A1").Value = "statusCode"
$ws3.Range("B1").Value = "IDC01S1"
$ws3.Range("C1").Value = "IDC01S2"
$ws3.Range("D1").Value = "IDC01S3"
$ws3.Range("E1").Value = "IDC01S4"

# deviceInstance row - numeric-looking text, stored with a leading
# apostrophe (quote-prefix) just like the source workbook.
$ws3.Range("A2").Value = "deviceInstance"
$ws3.Range("B2").Value = "'1"
$ws3.Range("C2").Value = "'2"
$ws3.Range("D2").Value = "'3"
$ws3.Range("E2").Value = "'4"

# deviceName row
$ws3.Range("A3").Value = "deviceName"
$ws3.Range("B3").Value = "Recycler1"
$ws3.Range("C3").Value = "Recycler2"
$ws3.Range("D3").Value = "Recycler3"
$ws3.Range("E3").Value = "Recycler4"

# deviceState row - also numeric-looking text with quote-prefix
$ws3.Range("A4").Value = "deviceState"
$ws3.Range("B4").Value = "'1"
$ws3.Range("C4").Value = "'2"
$ws3.Range("D4").Value = "'3"
$ws3.Range("E4").Value = "'4"

# Remaining descriptive rows - all placeholder "-"
$ws3.Range("A5").Value = "deviceDescription"
$ws3.Range("B5").Value = "-"
$ws3.Range("C5").Value = "-"
$ws3.Range("D5").Value = "-"
$ws3.Range("E5").Value = "-"

$ws3.Range("A6").Value = "deviceVendorIdentity"
$ws3.Range("B6").Value = "-"
$ws3.Range("C6").Value = "-"
$ws3.Range("D6").Value = "-"
$ws3.Range("E6").Value = "-"

$ws3.Range("A7").Value = "subDeviceName"
$ws3.Range("B7").Value = "-"
$ws3.Range("C7").Value = "-"
$ws3.Range("D7").Value = "-"
$ws3.Range("E7").Value = "-"

$ws3.Range("A8").Value = "subDeviceState"
$ws3.Range("B8").Value = "-"
$ws3.Range("C8").Value = "-"
$ws3.Range("D8").Value = "-"
$ws3.Range("E8").Value = "-"

$ws3.Range("A9").Value = "subDeviceDescription"
$ws3.Range("B9").Value = "-"
$ws3.Range("C9").Value = "-"
$ws3.Range("D9").Value = "-"
$ws3.Range("E9").Value = "-"

# Give Sheet3's columns a sensible best-fit-like width (A is a bit wider
# to fit the longest label, B:E share a common width).
$ws3.Cells.EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 3. Update selections / active sheet bookkeeping.
#    Sheet1 -> whole sheet selected (Ctrl+A style), no longer the active tab.
#    Sheet2 -> C1 selected.
#    Sheet3 -> becomes the active / tab-selected sheet (selected last).
# ---------------------------------------------------------------------
$ws1.Cells.Select() | Out-Null
$ws2.Range("C1").Select() | Out-Null
$ws3.Range("A1").Select() | Out-Null
